# Generate Report for Handback
#
# The "c4e3a62b-d9fe-4fda-8852-3a931081d1e1" file has finished its
# handback round-trip, so:
#   - its Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that tracks it
#     (Overview, zh-cn, de-de)
#   - the per-locale "Latest Handback DateTime" is stamped with the
#     handback completion time on the zh-cn and de-de detail sheets

$wb = $excel.ActiveWorkbook

$statusDone = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the c4e3a62b-... file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusDone   # zh-cn status
$overview.Range("C3").Value = $statusDone   # de-de status

# --- zh-cn detail sheet: row 3 is the c4e3a62b-... file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusDone
$zhcn.Range("H3").Value = "2016-03-22 12:13:23"

# --- de-de detail sheet: row 3 is the c4e3a62b-... file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusDone
$dede.Range("H3").Value = "2016-03-22 12:13:40"
